# Revert "fix[catalogue]: small bugfixes in ResourceDetails view"
#
# Net effect on this workbook:
#   1. Add an "admin" / "Manager" row to the "molgenis_members" sheet.
#   2. Make "molgenis_settings" the active/selected sheet again
#      (instead of "molgenis_members").

$wb = $excel.ActiveWorkbook

# 1) Append admin/Manager row to molgenis_members
$members = $wb.Worksheets.Item("molgenis_members")
$members.Range("A3").Value = "admin"
$members.Range("B3").Value = "Manager"

# 2) Re-activate molgenis_settings as the selected tab
$settings = $wb.Worksheets.Item("molgenis_settings")
$settings.Activate()
